$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: set text value "utd" (new shared string)
$ws.Range("C3").Value = "utd"

# C5: 76 -> 75
$ws.Range("C5").Value = 75

# C6: 630 -> 100
$ws.Range("C6").Value = 100

# C7: 643 -> 200
$ws.Range("C7").Value = 200

# D7: was text "+" -> now numeric 0
$ws.Range("D7").Value = 0

# Update selection to C3
$ws.Range("C3").Select()
